$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '63.952.47'
Set-TextValue $ws 'E2' '  -0.55%  '

Set-TextValue $ws 'D3' '2.640.22'
Set-TextValue $ws 'E3' '  -0.09%  '

Set-TextValue $ws 'E4' '  +0.00%  '

Set-TextValue $ws 'D5' '580.42'
Set-TextValue $ws 'E5' '  -0.08%  '

Set-TextValue $ws 'D6' '156.37'
Set-TextValue $ws 'E6' '  -0.20%  '

Set-TextValue $ws 'D7' '0.628'
Set-TextValue $ws 'E7' '  -3.61%  '

Set-TextValue $ws 'E8' '  +0.02%  '

Set-TextValue $ws 'D9' '2.636.99'
Set-TextValue $ws 'E9' '  -0.15%  '

Set-TextValue $ws 'E10' '  -3.50%  '

Set-TextValue $ws 'D11' '5.80'
Set-TextValue $ws 'E11' '  -0.50%  '

Set-TextValue $ws 'D12' '0.384'
Set-TextValue $ws 'E12' '  -1.40%  '

Set-TextValue $ws 'E13' '  +0.82%  '

Set-TextValue $ws 'D14' '28.65'
Set-TextValue $ws 'E14' '  -0.85%  '

Set-TextValue $ws 'D15' '3.112.74'
Set-TextValue $ws 'E15' '  -0.27%  '

Set-TextValue $ws 'E16' '  -1.53%  '

Set-TextValue $ws 'D17' '63.827.34'
Set-TextValue $ws 'E17' '  -0.48%  '

Set-TextValue $ws 'D18' '2.652.59'
Set-TextValue $ws 'E18' '  +0.63%  '

Set-TextValue $ws 'D19' '12.15'
Set-TextValue $ws 'E19' '  -0.96%  '

Set-TextValue $ws 'E20' '  +2.70%  '

Set-TextValue $ws 'E21' '  -3.81%  '

Set-TextValue $ws 'D22' '345.33'
Set-TextValue $ws 'E22' '  -0.85%  '

Set-TextValue $ws 'E23' '  +0.09%  '

Set-TextValue $ws 'D24' '68.16'
Set-TextValue $ws 'E24' '  +0.40%  '

Set-TextValue $ws 'E25' '  +5.69%  '

Set-TextValue $ws 'D26' '0.0000112'
Set-TextValue $ws 'E26' '  +1.69%  '

Set-TextValue $ws 'D27' '9.31'
Set-TextValue $ws 'E27' '  -1.12%  '

Set-TextValue $ws 'E28' '  +2.57%  '

Set-TextValue $ws 'D29' '584.39'
Set-TextValue $ws 'E29' '  -0.37%  '

Set-TextValue $ws 'D30' '8.17'
Set-TextValue $ws 'E30' '  +1.59%  '

Set-TextValue $ws 'B31' 'Binance-PegBSC-USD'
Set-TextValue $ws 'C31' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws 'D31' '1.00'
Set-TextValue $ws 'E31' '  -0.11%  '

Set-TextValue $ws 'B32' 'Kaspa'
Set-TextValue $ws 'C32' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws 'D32' '0.161'
Set-TextValue $ws 'E32' '  -0.65%  '

Set-TextValue $ws 'E33' '  -1.15%  '

Set-TextValue $ws 'E34' '  -0.07%  '

Set-TextValue $ws 'E35' '  +0.01%  '

Set-TextValue $ws 'D36' '5.46'
Set-TextValue $ws 'E36' '  +2.62%  '

Set-TextValue $ws 'D37' '0.403'
Set-TextValue $ws 'E37' '  -2.78%  '

Set-TextValue $ws 'D38' '19.76'
Set-TextValue $ws 'E38' '  -1.75%  '

Set-TextValue $ws 'E39' '  +0.02%  '

Set-TextValue $ws 'D40' '1.92'
Set-TextValue $ws 'E40' '  -0.37%  '

Set-TextValue $ws 'D41' '153.09'
Set-TextValue $ws 'E41' '  +0.60%  '

Set-TextValue $ws 'E42' '  -0.04%  '

Set-TextValue $ws 'E43' '  +6.80%  '

Set-TextValue $ws 'D44' '41.96'
Set-TextValue $ws 'E44' '  -0.26%  '

Set-TextValue $ws 'D45' '163.12'
Set-TextValue $ws 'E45' '  +2.14%  '

Set-TextValue $ws 'D46' '24.25'
Set-TextValue $ws 'E46' '  +3.33%  '

Set-TextValue $ws 'D47' '3.91'
Set-TextValue $ws 'E47' '  -2.98%  '

Set-TextValue $ws 'E48' '  -2.29%  '

Set-TextValue $ws 'D49' '0.635'
Set-TextValue $ws 'E49' '  -0.26%  '

Set-TextValue $ws 'E50' '  -3.28%  '

Set-TextValue $ws 'E51' '  -2.19%  '
